# Rename document/table attributes in the ObjTables header strings to
# lowerCamelCase (ObjTablesVersion -> objTablesVersion, Type -> type, Id -> id).

$wb = $excel.ActiveWorkbook

$tocSheet = $wb.Worksheets.Item("!!_Table of contents")
$tocSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$tocSheet.Range("A2").Value = "!!ObjTables type='TableOfContents'"

$schemaSheet = $wb.Worksheets.Item("!!Schema repo metadata")
$schemaSheet.Range("A1").Value = "!!ObjTables type='Data' id='SchemaRepoMetadata'"

$model1Sheet = $wb.Worksheets.Item("!!Model1s")
$model1Sheet.Range("A1").Value = "!!ObjTables type='Data' id='Model1'"
